$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D already has text-formatted cells (e.g. "43.792.03"); keep them as text
# by forcing the number format to Text before assigning values, so Excel does not
# auto-convert numeric-looking strings into real numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.757.48"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.312.62"
$ws.Range("E3").Value = "  +4.11%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "96.64"
$ws.Range("E5").Value = "  +4.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "269.96"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.621"
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.34"
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0944"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.00"
$ws.Range("E12").Value = "  -2.93%  "
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.655.55"
$ws.Range("E14").Value = "  +3.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.43"
$ws.Range("E15").Value = "  +2.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.865"
$ws.Range("E16").Value = "  +8.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.314.86"
$ws.Range("E17").Value = "  +3.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.709.73"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("E19").Value = "  +5.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.37"
$ws.Range("E20").Value = "  +5.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.56"
$ws.Range("E21").Value = "  +3.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "238.94"
$ws.Range("E22").Value = "  +2.95%  "
$ws.Range("E23").Value = "  -3.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.35"
$ws.Range("E24").Value = "  +3.84%  "
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.51"
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.30"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("E28").Value = "  -1.83%  "
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.35"
$ws.Range("E30").Value = "  +7.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.83"
$ws.Range("E31").Value = "  -8.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "174.37"
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0896"
$ws.Range("E33").Value = "  -2.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.45"
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("E35").Value = "  +2.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0360"
$ws.Range("E36").Value = "  +2.72%  "
$ws.Range("E37").Value = "  -3.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.35"
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("E39").Value = "  -7.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.243"
$ws.Range("E40").Value = "  +11.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.34"
$ws.Range("E41").Value = "  +8.92%  "
$ws.Range("E42").Value = "  +18.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.07"
$ws.Range("E43").Value = "  -3.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.14"
$ws.Range("E44").Value = "  +9.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.73"
$ws.Range("E45").Value = "  -2.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.32"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("E47").Value = "  +3.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "99.99"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.19"
$ws.Range("E49").Value = "  +0.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.190"
$ws.Range("E50").Value = "  +17.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.542.57"
$ws.Range("E51").Value = "  +4.12%  "
